# Generate Report for handback
# The localization "8320dd9d-..." item has completed its handoff/handback
# round-trip, so update the report rows that previously carried the
# "not yet processed" placeholder timestamp with the real handoff/handback
# datetimes, for both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-08 09:13:04"
$wsZhCn.Range("G3").Value = "2016-01-08 09:13:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-08 09:13:13"
$wsDeDe.Range("G3").Value = "2016-01-08 09:14:04"
